$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("program")

# --- Header row (row 1): replace the single "program_name" header with the
#     full set of REPROG_* / CED_* / BUSPAR_* column headers, columns A..N.

$headers = @(
    "REPROG_ID_PRE",
    "REPROG_TITLE",
    "CED_ID_PRE",
    "CED_NAME_PRE",
    "REPROG_ACTIVE_IND",
    "REPROG_COMMENT",
    "REPROG_UW_DEPARTMENT_CD",
    "REPROG_UW_DEPARTMENT_NAME",
    "REPROG_UW_DEPARTMENT_LOB_CD",
    "REPROG_UW_DEPARTMENT_LOB_NAME",
    "BUSPAR_CED_REG_CLASS_CD",
    "BUSPAR_CED_REG_CLASS_NAME",
    "REPROG_MAIN_CURRENCY_CD",
    "REPROG_MANAGEMENT_REPORTING_LOB_CD"
)

# A1 already carries the bold / centered / bordered header style. Copy that
# formatting across the rest of row 1 (B1:N1) before writing the header
# text, so every header cell ends up with the same style.
$ws.Cells.Item(1, 1).Copy($ws.Range("B1:N1"))

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}

# --- Data row (row 2): program id, program title (the previous sole
#     value), the active flag, and blank placeholders for every other
#     column of the new data model. ---

$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "AVIATION_AXA_XL_2024"
$ws.Cells.Item(2, 5).Value = $true

# The remaining data-model columns have no value yet, but the row still
# carries a (blank) text cell for each of them. A bare "" assignment
# deletes a cell outright, so force an empty text entry the way Excel's
# UI does (leading apostrophe = enter as text), then drop back to the
# Normal style so no stray quote-prefix formatting is left behind.
$blankColumns = @(3, 4, 6, 7, 8, 9, 10, 11, 12, 13, 14)
foreach ($col in $blankColumns) {
    $cell = $ws.Cells.Item(2, $col)
    $cell.Value = "'"
    $cell.Style = "Normal"
}
